# Update donor / donation record data on Sheet1, row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = "Anne"              # Donor First Name
$ws.Range("G2").Value = "Donk"              # Donor Last Name
$ws.Range("M2").Value = "annab@gmail.com"   # Donor Email
$ws.Range("O2").Value = "tester26"          # Reference code
$ws.Range("A2").Value = "AB522581053"       # Receipt ID
